$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 112231329
$ws.Range("B2").Value = 77719
$ws.Range("E2").Value = 283
$ws.Range("F2").Value = "Kavernularia"
$ws.Range("G2").Value = "Hypogymnia hultenii"
$ws.Range("H2").Value = "(Degel.) Krog"

$ws.Range("A3").Value = 112231350
$ws.Range("B3").Value = 73820
$ws.Range("E3").Value = 6440
$ws.Range("F3").Value = "Vitgrynig nållav"
$ws.Range("G3").Value = "Chaenotheca subroscida"
$ws.Range("H3").Value = "(Eitner) Zahlbr."

$ws.Range("B4").Value = 85836

$ws.Range("A5").Value = 112231361
$ws.Range("B5").Value = 77636
$ws.Range("E5").Value = 6425
$ws.Range("F5").Value = "Garnlav"
$ws.Range("G5").Value = "Alectoria sarmentosa"
$ws.Range("H5").Value = "(Ach.) Ach."

$ws.Range("A6").Value = 112231346
$ws.Range("B6").Value = 76620
$ws.Range("E6").Value = 228579
$ws.Range("F6").Value = "Liten svartspik"
$ws.Range("G6").Value = "Chaenothecopsis nana"
$ws.Range("H6").Value = "Tibell"
